$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToPutIntoProg")

# Replace the hard-coded placeholder values (1234) in column C (rows 15-180)
# with a lookup formula that pulls the correct exchange-rate value from the
# K/L reference table (K15:K187 / L15:L187) based on the currency code in
# column A. Applying the formula in the same chunks that already exist for
# the neighbouring column E reproduces the same shared-formula grouping.
$ws.Range("C15").Formula = '=INDEX($K$15:$K$187,MATCH(A15,$L$15:$L$187,0))'
$ws.Range("C16:C79").Formula = '=INDEX($K$15:$K$187,MATCH(A16,$L$15:$L$187,0))'
$ws.Range("C80:C143").Formula = '=INDEX($K$15:$K$187,MATCH(A80,$L$15:$L$187,0))'
$ws.Range("C144:C180").Formula = '=INDEX($K$15:$K$187,MATCH(A144,$L$15:$L$187,0))'

# Recalculate so the dependent column E (which concatenates the new C values
# into the generated Java source lines) picks up the refreshed numbers.
$excel.CalculateFull()

# Update the sheet's active cell / selection from N11 to H11.
$ws.Activate() | Out-Null
$ws.Range("H11").Select() | Out-Null
